$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:E51 as Text so numeric-looking price/volume strings are preserved
# exactly as text (matching the source data which stores them as inline strings)
# rather than being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.238.15'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '1.805.02'
$ws.Range("E3").Value = '  +2.47%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '339.11'
$ws.Range("E5").Value = '  +1.00%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").Value = '0.4912'
$ws.Range("E7").Value = '  +29.86%  '

$ws.Range("D8").Value = '0.3771'
$ws.Range("E8").Value = '  +12.09%  '

$ws.Range("D9").Value = '45.55'
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").Value = '0.07739'
$ws.Range("E10").Value = '  +7.56%  '

$ws.Range("D11").Value = '1.149'
$ws.Range("E11").Value = '  +2.14%  '

$ws.Range("D12").Value = '22.64'
$ws.Range("E12").Value = '  +0.67%  '

$ws.Range("D13").Value = '0.9999'
$ws.Range("E13").Value = '  +0.01%  '

$ws.Range("D14").Value = '6.347'
$ws.Range("E14").Value = '  +2.28%  '

$ws.Range("D15").Value = '7.332'
$ws.Range("E15").Value = '  +1.78%  '

$ws.Range("D16").Value = '1.803.88'
$ws.Range("E16").Value = '  +2.61%  '

$ws.Range("D17").Value = '0.00001100'
$ws.Range("E17").Value = '  +4.16%  '

$ws.Range("D18").Value = '0.06744'
$ws.Range("E18").Value = '  +2.33%  '

$ws.Range("D19").Value = '82.31'
$ws.Range("E19").Value = '  +2.24%  '

$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = '17.47'
$ws.Range("E21").Value = '  +2.60%  '

$ws.Range("D22").Value = '6.450'
$ws.Range("E22").Value = '  +2.61%  '

$ws.Range("D23").Value = '28.188.77'
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").Value = '11.99'
$ws.Range("E24").Value = '  +2.21%  '

$ws.Range("D25").Value = '2.396'
$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("D26").Value = '20.98'
$ws.Range("E26").Value = '  +5.63%  '

$ws.Range("D27").Value = '2.421'
$ws.Range("E27").Value = '  +3.27%  '

$ws.Range("D28").Value = '151.71'
$ws.Range("E28").Value = '  -1.11%  '

$ws.Range("D29").Value = '2.007.23'
$ws.Range("E29").Value = '  +2.38%  '

$ws.Range("D30").Value = '134.40'
$ws.Range("E30").Value = '  +2.25%  '

$ws.Range("D31").Value = '1.278'
$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("D32").Value = '4.047'
$ws.Range("E32").Value = '  +0.82%  '

$ws.Range("D33").Value = '0.09882'
$ws.Range("E33").Value = '  +12.09%  '

$ws.Range("D34").Value = '5.960'
$ws.Range("E34").Value = '  +2.46%  '

$ws.Range("D35").Value = '0.02403'
$ws.Range("E35").Value = '  +2.21%  '

$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").Value = '0.2241'
$ws.Range("E36").Value = '  +5.79%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '12.27'
$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = '0.06390'
$ws.Range("E38").Value = '  +3.05%  '

$ws.Range("D39").Value = '0.6736'
$ws.Range("E39").Value = '  +1.60%  '

$ws.Range("D40").Value = '5.247'
$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("E41").Value = '  +2.15%  '

$ws.Range("E42").Value = '  +0.40%  '

$ws.Range("D43").Value = '8.160'
$ws.Range("E43").Value = '  +1.26%  '

$ws.Range("D44").Value = '14.14'
$ws.Range("E44").Value = '  +2.52%  '

$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").Value = '0.6196'
$ws.Range("E46").Value = '  +2.31%  '

$ws.Range("D47").Value = '3.880'
$ws.Range("E47").Value = '  +1.70%  '

$ws.Range("D48").Value = '129.53'
$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("D49").Value = '2.064'
$ws.Range("E49").Value = '  +2.29%  '

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.175'
$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.07121'
$ws.Range("E51").Value = '  -1.13%  '
